$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Configuration")
$ws.Range("C12").Value = "http://localhost:8080"
$ws.Range("C12").Select()
